$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row (A1:K1) is now bold (it already had centered alignment from
# style index 1; Excel creates a new font + cellXf combining bold + center).
$ws.Range("A1:K1").Font.Bold = $true

# The "target" column (I2:I10) values were relabeled from "proton" to "p".
# Writing "p" to every data row leaves "proton" with no remaining
# references, so it is dropped from the shared-string table on save (and
# the new string "p" is appended) -- matching the diff exactly.
$ws.Range("I2:I10").Value = "p"

# The sheet's saved selection moved from L7 to the column that was just
# edited (I2:I10, active cell I2).
[void]$ws.Range("I2:I10").Select()
